$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (QBTS)
$ws.Range("D2").Value = 27
$ws.Range("E2").Value = 59.2
$ws.Range("F2").Value = 19.1
$ws.Range("K2").Value = 61.7
$ws.Range("N2").Value = 52.28493729186943

# Row 3 (IONQ)
$ws.Range("D3").Value = 52.69
$ws.Range("E3").Value = 58.9
$ws.Range("F3").Value = 6.88
$ws.Range("H3").Value = 56
$ws.Range("K3").Value = 57.7
$ws.Range("N3").Value = 52.28493729186943

# Row 4 (RGTI)
$ws.Range("D4").Value = 28.11
$ws.Range("E4").Value = 56.5
$ws.Range("F4").Value = 9.93
$ws.Range("K4").Value = 57.1
$ws.Range("N4").Value = 52.28493729186943

# Row 5 (IBM)
$ws.Range("D5").Value = 307.94
$ws.Range("E5").Value = 52.1
$ws.Range("F5").Value = -0.21
$ws.Range("G5").Value = 40
$ws.Range("K5").Value = 54.1
$ws.Range("N5").Value = 52.28493729186943
